$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Bump the "GroupDocs.Assembly" evaluation-watermark version that appears
#    in the first paragraph of the document body.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "GroupDocs.Assembly 25.6.", $true, $false, $false, $false, $false,
    $true, 1, $false, "GroupDocs.Assembly 25.12.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Register the built-in "Hyperlink" character style that newer
#    GroupDocs.Assembly / Aspose.Words releases emit in the style sheet
#    (basedOn DefaultParagraphFont, visually blue + underlined).
# ---------------------------------------------------------------------------
$hyperlinkStyle = $d.Styles.Add("Hyperlink", 2)
$hyperlinkStyle.BaseStyle = $d.Styles("DefaultParagraphFont")
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.Font.Color = 12673797
$hyperlinkStyle.Font.Underline = 1
